# Edit script: add "Verification Feedback" row to Summary sheet and
# add a new "Content Verification" worksheet with fact-check data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Update the "Summary" sheet: insert a new row 7 "Verification
#    Feedback" above the existing "Generated On" row, and refresh
#    the "Generated On" timestamp.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Summary")
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "Verification Feedback"
$ws.Range("B7").Value = "Overall Evaluation Summary: The marking scheme is factually accurate and aligns with the official terminology, institutional roles, and educational pathways defined by the VTC. The marking breakdowns are logical and prioritize the most critical information. Final Recommendations: 1. Terminology Consistency: Ensure that 'Vocational and Professional Education and Training' is always written in full for the 2-mark component of Q1 to avoid ambiguity. 2. Partial Marks: The 'General Grading Guide' is helpful, but for Q2 and Q4, ensure examiners know that naming the specific qualification (Higher Diploma vs. Degree or DFS) is the 'hurdle' requirement for moving above the 5-mark threshold. 3. Wording Improvement: In Q4, the phrase 'guarantees entry' should be used cautiously; 'eligibility to apply' or 'articulation pathway' is more technically accurate as some HD programs have specific interviews or subject requirements. Citations: vtc.edu.hk, iabhongkong.com, legco.gov.hk."

$ws.Range("B8").Value = "2026-01-07 13:54:45"

# ---------------------------------------------------------------
# 2. Add a new "Content Verification" worksheet at the end of the
#    workbook containing the per-question fact-check results.
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cv = $wb.Worksheets.Add($null, $lastSheet)
$cv.Name = "Content Verification"

$cv.Range("A1").Value = "question_number"
$cv.Range("B1").Value = "is_correct"
$cv.Range("C1").Value = "feedback"
$cv.Range("D1").Value = "suggestion"
$cv.Range("A2").Value = "Q1"
$cv.Range("B2").Value = $true
$cv.Range("C2").Value = "Fact Check: Correct. VPET stands for Vocational and Professional Education and Training. The VTC is officially recognized as the largest provider of VPET in Hong Kong. Verification Detail: Official VTC corporate information and annual reports confirm this acronym and the council's status as the primary statutory body for vocational training."
$cv.Range("D2").Value = "The answer is solid. To make the 'Importance' section even more robust, you could mention that VPET is a 'valued choice' alongside traditional academic paths, which is a key part of the Hong Kong government's recent branding of the sector."
$cv.Range("A3").Value = "Q2"
$cv.Range("B3").Value = $true
$cv.Range("C3").Value = "Fact Check: Correct. IVE (Hong Kong Institute of Vocational Education) is the primary provider of Higher Diplomas (HD) and Certificates. THEi (Technological and Higher Education Institute of Hong Kong) was specifically established to offer vocationally-oriented Bachelor’s Degrees. Verification Detail: While THEi does offer a small number of Higher Diplomas, its defining characteristic in the VTC ecosystem is providing the degree-level 'top' of the progression ladder."
$cv.Range("D3").Value = "The marking breakdown is fair. You might add a note for examiners that mentioning 'Applied Degrees' is also a highly accurate description for THEi."
$cv.Range("A4").Value = "Q3"
$cv.Range("B4").Value = $true
$cv.Range("C4").Value = "Fact Check: Correct. 'Think and Do' is the official motto and pedagogical approach of the VTC. Verification Detail: VTC branding materials emphasize this approach to distinguish their training from purely theoretical academic studies. It specifically refers to the integration of 'Head' (knowledge) and 'Hands' (skills)."
$cv.Range("D4").Value = "The marking breakdown is excellent. Ensure students understand that 'Think' refers to professional knowledge/theory, not just general thinking."
$cv.Range("A5").Value = "Q4"
$cv.Range("B5").Value = $true
$cv.Range("C5").Value = "Fact Check: Correct. The Diploma of Foundation Studies (DFS) is the standard one-year program for S6 leavers to bridge into a Higher Diploma. The Diploma of Vocational Education (DVE) is also a valid pathway (offered by Youth College). Verification Detail: VTC admission guides list DFS as the primary 'Level 3' qualification on the Qualifications Framework (QF) that grants eligibility for HD programs."
$cv.Range("D5").Value = "You may want to mention the Diploma of Applied Education (DAE) as an alternative, though DFS/DVE are the specific VTC-branded programs. The current answer is the most relevant for a VTC-specific context."
$cv.Range("A6").Value = "Q5"
$cv.Range("B6").Value = $true
$cv.Range("C6").Value = "Fact Check: Correct. VTC operates 25 Training Boards that represent various industries to ensure curriculum relevance. Verification Detail: Key benefits include Work-Integrated Learning (WIL), which is the formal term VTC uses for internships, and the use of industry-donated facilities (e.g., the Boeing-standard engine in their engineering labs)."
$cv.Range("D6").Value = "The examples provided (Internships, Job Placement, Equipment) are the most common and correct answers. Using the term 'Work-Integrated Learning (WIL)' could be a 'bonus' keyword for a 10-mark answer."

# Copy the bold/bordered header style used elsewhere in the workbook
# (e.g. Validation!A1:B1) onto the new header row.
$headerSrc = $wb.Worksheets.Item("Validation")
$headerSrc.Range("A1:B1").Copy()
$cv.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$cv.Range("A1").Select()
